# Updated symbol list on Thu Feb 16 21:45:48 UTC 2023 with GitHub Actions
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) values for the
# crypto coin rows on the active sheet, matching the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ column letter = new text value }
# Only the cells whose Price/Volume text actually changed are listed here.
$updates = [ordered]@{
    2  = @{ D = "312.33";    E = "-0.02%" }
    3  = @{ D = "48.71";     E = "9.62%" }
    4  = @{ D = "5.266";     E = "2.91%" }
    5  = @{ D = "0.07877";   E = "-1.55%" }
    6  = @{ D = "4.575";     E = "1.72%" }
    7  = @{ D = "1.322";     E = "21.93%" }
    8  = @{ D = "1.614";     E = "-2.60%" }
    9  = @{ D = "0.1239";    E = "-4.41%" }
    10 = @{ D = "0.1947";    E = "2.13%" }
    11 = @{ D = "0.09476";   E = "0.94%" }
    12 = @{ D = "0.04533";   E = "7.35%" }
    13 = @{ D = "0.1049";    E = "0.93%" }
    14 = @{ D = "0.001309";  E = "0.31%" }
    15 = @{            E = "0.32%" }
    16 = @{ D = "0.005887";  E = "0.85%" }
    17 = @{ D = "3.344";     E = "-0.97%" }
    18 = @{            E = "2.72%" }
    19 = @{            E = "1.73%" }
    20 = @{            E = "-0.05%" }
    21 = @{ D = "0.1366";    E = "-0.25%" }
    22 = @{ D = "0.3068";    E = "-1.97%" }
    23 = @{ D = "0.001294";  E = "1.54%" }
    24 = @{ D = "0.004178";  E = "-9.13%" }
    25 = @{            E = "1.12%" }
    26 = @{ D = "0.0003554" }
    38 = @{ D = "0.02630";   E = "-0.84%" }
    39 = @{ D = "0.05863";   E = "8.56%" }
    40 = @{            E = "91.96%" }
    41 = @{ D = "0.008021";  E = "2.76%" }
    42 = @{ D = "0.1442";    E = "2.16%" }
    43 = @{ D = "0.008301";  E = "13.23%" }
    44 = @{ D = "0.008660";  E = "9.86%" }
    45 = @{ D = "0.3142";    E = "0.86%" }
    46 = @{ D = "0.00006952"; E = "3.31%" }
    47 = @{            E = "1.18%" }
    49 = @{ D = "0.004016";  E = "1.18%" }
    50 = @{ D = "0.00002109"; E = "1.18%" }
    51 = @{ D = "0.0002008"; E = "1.18%" }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $cellRef = "$col$row"
        $newValue = $updates[$row][$col]
        $cell = $ws.Range($cellRef)
        # Force the cell to keep storing plain text (matching the workbook's
        # existing inline-string cells) instead of letting Excel reinterpret
        # numeric-looking / percent-looking text as a number.
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
    }
}
